# Weekly update: insert two new rows of "Camote" price data (Vega Central
# Mapocho de Santiago) at the top of the existing data block, pushing the
# older rows down. This mirrors a new week's worth of readings being
# prepended above the historical rows that already occupied 126-140.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 126; everything currently at 126 and below
# (through 140) shifts down to 128-142.
$ws.Rows("126:127").Insert()

# New row 126
$ws.Range('A126').Value = 9
$ws.Range('B126').Value = 'Vega Central Mapocho de Santiago'
$ws.Range('C126').Value = 'Metropolitana'
$ws.Range('D126').Value = 44984
$ws.Range('E126').Value = 13
$ws.Range('F126').Value = 100114002
$ws.Range('G126').Value = 'Camote'
$ws.Range('H126').Value = 'Sin especificar'
$ws.Range('I126').Value = 'Primera'
$ws.Range('J126').Value = 790
$ws.Range('K126').Value = 18000
$ws.Range('L126').Value = 19000
$ws.Range('M126').Value = 18494
$ws.Range('N126').Value = '$/caja 18 kilos'
$ws.Range('O126').Value = 'Perú'
$ws.Range('P126').Value = 1027
$ws.Range('Q126').Value = 18
$ws.Range('R126').Value = 'Hortaliza'

# New row 127
$ws.Range('A127').Value = 9
$ws.Range('B127').Value = 'Vega Central Mapocho de Santiago'
$ws.Range('C127').Value = 'Metropolitana'
$ws.Range('D127').Value = 44984
$ws.Range('E127').Value = 13
$ws.Range('F127').Value = 100114002
$ws.Range('G127').Value = 'Camote'
$ws.Range('H127').Value = 'Sin especificar'
$ws.Range('I127').Value = 'Primera'
$ws.Range('J127').Value = 970
$ws.Range('K127').Value = 15000
$ws.Range('L127').Value = 16000
$ws.Range('M127').Value = 15546
$ws.Range('N127').Value = '$/malla 18 kilos'
$ws.Range('O127').Value = 'Perú'
$ws.Range('P127').Value = 864
$ws.Range('Q127').Value = 18
$ws.Range('R127').Value = 'Hortaliza'
